$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 68

# Column A: copy the format of the cell above (style index "1"), then set its value
$ws.Cells.Item($newRow - 1, 1).Copy($ws.Cells.Item($newRow, 1))
$ws.Cells.Item($newRow, 1).Value = 66

# Column B: enter the date as a text literal (via a text formula) so Excel does not
# auto-convert it to a date serial number, then collapse the formula down to a plain
# value in place so it is stored as a shared string, matching the other date cells
# (which carry no explicit style).
$dateCell = $ws.Cells.Item($newRow, 2)
$dateCell.Formula = "=""06-09-2020"""
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)

$ws.Cells.Item($newRow, 3).Value = 15238
$ws.Cells.Item($newRow, 4).Value = 262
$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 14976
